$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date for every existing data row (2-23)
# from 45192 (2023-09-23) to 45202 (2023-10-03).
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# Row 23 gains an explicit row height (matches the rest of the sheet's rows).
$ws.Rows.Item(23).RowHeight = 15

# Append a new record in row 24.
$ws.Cells.Item(24, 1).Value = "A 47107-2023"
$ws.Cells.Item(24, 2).Value = 45196
$ws.Cells.Item(24, 3).Value = 45202
$ws.Cells.Item(24, 4).Value = "OKÄNT"
$ws.Cells.Item(24, 5).Value = "OKÄNT"
$ws.Cells.Item(24, 7).Value = 1.5
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 0
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).Value = 0
$ws.Cells.Item(24, 14).Value = 0
$ws.Cells.Item(24, 15).Value = 0
$ws.Cells.Item(24, 16).Value = 0
$ws.Cells.Item(24, 17).Value = 0

# B24/C24 (the two date cells) use the same YYYY-MM-DD date format as the rest
# of the sheet's B/C columns.
$ws.Range("B24:C24").NumberFormat = "YYYY-MM-DD"

# R24 carries no value, just the wrap-text style used by the rest of column R.
$ws.Range("R24").WrapText = $true
